$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Good Morning" -> "GIT UPDATE" (the jgit commit message's text update)
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the new active cell/selection recorded in the saved view
$ws.Range("E8").Select()
